$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 100.75
$ws.Range("I4").Value = 100.75
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 100.75
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 13.25

$ws.Range("H18").Value = 749.4
$ws.Range("I18").Value = 749.4
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 749.4
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -465.4

$ws.Range("H28").Value = 1109.7142
$ws.Range("I28").Value = 449.75
$ws.Range("J28").Value = 1989.6666
$ws.Range("K28").Value = 449.75
$ws.Range("L28").Value = 1989.6666
$ws.Range("M28").Value = 35.25
$ws.Range("N28").Value = -2959.6666

$ws.Range("H31").Value = 306.57144
$ws.Range("I31").Value = 306.57144
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 919.71432
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -689.71432
$ws.Range("N31").ClearContents()

$ws.Range("H51").Value = 6998
$ws.Range("I51").Value = 6998
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 6998
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -6514

$ws.Range("H106").Value = 1002.5714
$ws.Range("I106").Value = 836.3333
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 836.3333
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -205.3333
$ws.Range("N106").Value = -3262

$ws.Range("H132").Value = 2743.7778
$ws.Range("I132").Value = 2737
$ws.Range("J132").Value = 2798
$ws.Range("K132").Value = 8211
$ws.Range("L132").Value = 8394
$ws.Range("M132").Value = -5681
$ws.Range("N132").Value = -13454

$ws.Range("H135").Value = 400
$ws.Range("I135").Value = 400
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3600
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1065

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 4223.1113
$ws.Range("I26").Value = 3000
$ws.Range("J26").Value = 4376
$ws.Range("K26").Value = 3000
$ws.Range("L26").Value = 4376
$ws.Range("M26").Value = -2670
$ws.Range("N26").Value = -5036

$ws.Range("H74").Value = 447856.78
$ws.Range("I74").Value = 503213.88
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 503213.88
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -502339.88
$ws.Range("N74").Value = -6748

$ws.Range("H77").Value = 447856.78
$ws.Range("I77").Value = 503213.88
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 2516069.4
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -2511701.4
$ws.Range("N77").Value = -33736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 590
$ws.Range("I22").Value = 510
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 510
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -337
$ws.Range("N22").Value = -1096

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H99").Value = 2772.5
$ws.Range("I99").Value = 2772.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2772.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1274.5

$ws.Range("H134").Value = 4729
$ws.Range("I134").Value = 4729
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14187
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -11652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1338.0952
$ws.Range("I15").Value = 655
$ws.Range("J15").Value = 15000
$ws.Range("K15").Value = 655
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = -485
$ws.Range("N15").Value = -15340

$ws.Range("H35").Value = 4008.6667
$ws.Range("I35").Value = 3804.8
$ws.Range("J35").Value = 5028
$ws.Range("K35").Value = 3804.8
$ws.Range("L35").Value = 5028
$ws.Range("M35").Value = -3510.8
$ws.Range("N35").Value = -5616

$ws.Range("H55").Value = 40666.668
$ws.Range("I55").Value = 41000
$ws.Range("J55").Value = 40000
$ws.Range("K55").Value = 41000
$ws.Range("L55").Value = 40000
$ws.Range("M55").Value = -40685
$ws.Range("N55").Value = -40630

$ws.Range("H132").Value = 4418.8823
$ws.Range("I132").Value = 3638.4167
$ws.Range("J132").Value = 6292
$ws.Range("K132").Value = 10915.2501
$ws.Range("L132").Value = 18876
$ws.Range("M132").Value = -8385.250100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 75745.03999999999
$ws.Range("I4").Value = 40265.25
$ws.Range("J4").Value = 501502.5
$ws.Range("K4").Value = 120795.75
$ws.Range("L4").Value = 1504507.5
$ws.Range("M4").Value = -120683.75
$ws.Range("N4").Value = -1504731.5

$ws.Range("H6").Value = 36
$ws.Range("I6").Value = 47.333332
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 141.999996
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = -28.99999600000001
$ws.Range("N6").Value = -232

$ws.Range("H26").Value = 323.6
$ws.Range("I26").Value = 127.28571
$ws.Range("J26").Value = 781.6667
$ws.Range("K26").Value = 381.85713
$ws.Range("L26").Value = 2345.0001
$ws.Range("M26").Value = -93.85712999999998
$ws.Range("N26").Value = -2921.0001

$ws.Range("H46").Value = 4556.4
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 5495.5
$ws.Range("K46").Value = 2400
$ws.Range("L46").Value = 16486.5
$ws.Range("M46").Value = -2309
$ws.Range("N46").Value = -16668.5

$ws.Range("H114").Value = 3500
$ws.Range("I114").Value = 4166.6665
$ws.Range("J114").Value = 1500
$ws.Range("K114").Value = 12499.9995
$ws.Range("L114").Value = 4500
$ws.Range("M114").Value = -9245.999500000002
$ws.Range("N114").Value = -11008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1582
$ws.Range("I5").Value = 1582
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1582
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1470

$ws.Range("H24").Value = 18999
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 18999
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 18999
$ws.Range("N24").Value = -19345

$ws.Range("H132").Value = 4187.25
$ws.Range("I132").Value = 3374.75
$ws.Range("J132").Value = 4999.75
$ws.Range("K132").Value = 10124.25
$ws.Range("L132").Value = 14999.25
$ws.Range("M132").Value = -7594.25
$ws.Range("N132").Value = -20059.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 105.333336
$ws.Range("I2").Value = 105.42857
$ws.Range("J2").Value = 104
$ws.Range("K2").Value = 105.42857
$ws.Range("L2").Value = 104
$ws.Range("M2").Value = 6.571430000000007
$ws.Range("N2").Value = -328

$ws.Range("H46").Value = 700
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 700
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 700
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -1076

$ws.Range("H55").Value = 852.4
$ws.Range("I55").Value = 815.75
$ws.Range("J55").Value = 999
$ws.Range("K55").Value = 815.75
$ws.Range("L55").Value = 999
$ws.Range("M55").Value = -642.75
$ws.Range("N55").Value = -1345

$ws.Range("H100").Value = 4166.6665
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 4750
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 4750
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -5832

$ws.Range("H122").Value = 3612.375
$ws.Range("I122").Value = 3266.6667
$ws.Range("J122").Value = 3819.8
$ws.Range("K122").Value = 9800.000100000001
$ws.Range("L122").Value = 11459.4
$ws.Range("M122").Value = -7350.000100000001
$ws.Range("N122").Value = -16359.4

$ws.Range("H136").Value = 867833.8
$ws.Range("I136").Value = 1274750.8
$ws.Range("J136").Value = 54000
$ws.Range("K136").Value = 3824252.4
$ws.Range("L136").Value = 162000
$ws.Range("M136").Value = -3821702.4
$ws.Range("N136").Value = -167100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1500
$ws.Range("I81").Value = 1500
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1939
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 1500
$ws.Range("I84").Value = 1500
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -9696
$ws.Range("N84").ClearContents()

$ws.Range("H136").Value = 1015.7931
$ws.Range("I136").Value = 957.1667
$ws.Range("J136").Value = 1297.2
$ws.Range("K136").Value = 2871.5001
$ws.Range("L136").Value = 3891.6
$ws.Range("M136").Value = -321.5001000000002
$ws.Range("N136").Value = -8991.6

$ws.Range("H139").Value = 20650
$ws.Range("I139").Value = 20650
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 20650
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -15510
